$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B17 first: change from numeric 2 to text "2" (same displayed value, stored as text)
# This gets inserted into the shared string table before the B8 value, to match
# the shared string ordering produced by the original edit.
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "2"

# Update B8: change from numeric 8535 to text "6000"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "6000"

# Update the active selection to B8 (matches the saved selection in the diff)
$ws.Range("B8").Select()
